# Bump the "Förändrad" (Changed) date column (C) by one day for every data row.
# Column C contains a date serial number (e.g. 46074 -> 2026-02-21) that is
# updated to 46075 (2026-02-22) for all rows from row 2 through row 499.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 499
$range = $ws.Range("C2:C$lastRow")
$range.Value = 46075
